$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.318.36'
$ws.Range('E2').Value = '  -1.54%  '

$ws.Range('D3').Value = '3.316.76'
$ws.Range('E3').Value = '  -2.60%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.79'
$ws.Range('E5').Value = '  -3.59%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '648.04'
$ws.Range('E6').Value = '  -1.49%  '

$ws.Range('E7').Value = '  -7.94%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.414'
$ws.Range('E8').Value = '  -3.20%  '

$ws.Range('E9').Value = '  +0.16%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.973'
$ws.Range('E10').Value = '  -7.46%  '

$ws.Range('D11').Value = '3.315.04'
$ws.Range('E11').Value = '  -2.22%  '

$ws.Range('E12').Value = '  -4.45%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '39.61'
$ws.Range('E13').Value = '  -4.88%  '

$ws.Range('D14').Value = '96.061.41'
$ws.Range('E14').Value = '  -1.41%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.07'
$ws.Range('E15').Value = '  -3.11%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000249'
$ws.Range('E16').Value = '  -3.94%  '

$ws.Range('D17').Value = '3.935.40'
$ws.Range('E17').Value = '  -2.56%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.41'
$ws.Range('E18').Value = '  -1.03%  '

$ws.Range('D19').Value = '3.314.75'
$ws.Range('E19').Value = '  -2.34%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.63'
$ws.Range('E20').Value = '  -4.13%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.469'
$ws.Range('E21').Value = '  -2.98%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '494.68'
$ws.Range('E22').Value = '  -2.62%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.32'
$ws.Range('E23').Value = '  -4.33%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.30'
$ws.Range('E24').Value = '  -4.53%  '

$ws.Range('E25').Value = '  -5.39%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.39'
$ws.Range('E26').Value = '  +4.10%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '91.77'
$ws.Range('E27').Value = '  -5.02%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.78'
$ws.Range('E28').Value = '  -6.17%  '

$ws.Range('D29').Value = '3.491.01'
$ws.Range('E29').Value = '  -2.49%  '

$ws.Range('E30').Value = '  +0.43%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.140'
$ws.Range('E31').Value = '  -6.86%  '

$ws.Range('E32').Value = '  -7.69%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.184'
$ws.Range('E33').Value = '  -4.93%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.44'
$ws.Range('E34').Value = '  +12.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.60%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.538'
$ws.Range('E36').Value = '  -4.30%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '27.73'
$ws.Range('E37').Value = '  -6.27%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.48'
$ws.Range('E38').Value = '  +7.22%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.43'
$ws.Range('E39').Value = '  -4.55%  '

$ws.Range('E40').Value = '  -0.02%  '

$ws.Range('E41').Value = '  -4.58%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '500.17'
$ws.Range('E42').Value = '  -2.48%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '24.45'
$ws.Range('E43').Value = '  -1.01%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.64'
$ws.Range('E44').Value = '  -1.22%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.817'
$ws.Range('E45').Value = '  -3.40%  '

$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0400'
$ws.Range('E46').Value = '  -7.37%  '

$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.26'
$ws.Range('E47').Value = '  +0.56%  '

$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.38'
$ws.Range('E48').Value = '  +1.87%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.60'
$ws.Range('E49').Value = '  +1.79%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.53'
$ws.Range('E50').Value = '  +3.64%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.10'
$ws.Range('E51').Value = '  -6.21%  '
